$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "position" column (B) should now uniformly report -1 instead of an
# incrementing per-attribute position, since position is no longer tracked
# by the metadata generation step.
$ws.Range("B2:B7").Value = -1
